$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: change style from (4/5) to (6/7), matching the alternating-shade row pattern used elsewhere
$ws.Range("A5:E5").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)

# Row 36: new dialogue line, copy format from row 6 (style 4/5, ht 43.2), includes column A (filename)
$ws.Range("A6:E6").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)

# Row 37: second line of same dialogue, copy format from row 7 (style 4/5, no column A)
$ws.Range("B7:E7").Copy()
$ws.Range("B37:E37").PasteSpecial(-4122)

$ws.Range("A36").Value = "SCRIPT/T01P01A/us2305.ssb"
$ws.Range("B36").Value = 19
$ws.Range("B37").Value = 22

$ws.Range("C36").Value = ' Oooh, I must have it! The [CS:I]Lost\nLoot[CR] calls to me!'
$ws.Range("C37").Value = ' The [CS:I]Lost Loot[CR] is rightfully mine!'

$ws.Range("D36").Value = ' Оооо, мне оно нужно! [CS:I]Утерянное\nДобро[CR] взывает ко мне!'
$ws.Range("D37").Value = ' [CS:I]Утерянное Добро[CR] только моё!'

$ws.Range("E36").Value = ' Ïïïï, íîå ïîï îôçîï! [CS:I]Ôóåñÿîîïå\nÄïáñï[CR] âèúâàåó ëï íîå!'
$ws.Range("E37").Value = ' [CS:I]Ôóåñÿîîïå Äïáñï[CR] óïìûëï íïæ!'

$ws.Rows.Item(36).RowHeight = 43.2
$ws.Rows.Item(37).RowHeight = 27

$ws.Range("B41").Select()

Write-Output "done"
